$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 18 (Web and Network Science, assignment 2, due 24-Mar-2023, status "No") ---
# First copy row 17's current formatting (red "Not submitted" style) onto row 18,
# since that is the style the new row should end up with.
$ws.Range("A17:D17").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A18").Value = "Web and Network Science"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 45009
$ws.Range("D18").Value = "No"

# --- Row 17: Embedded Image Processing assignment 4 moves from "Not submitted" (red) to "In progress" (orange) ---
# Copy row 12's current formatting (orange "in progress" style) onto row 17.
$ws.Range("A12:D12").Copy() | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 12: Embedded Image Processing assignment 3 moves from "In progress" (orange) to "Submitted" (green) ---
# Copy row 10's formatting (green "submitted" style) onto row 12.
$ws.Range("A10:D10").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E18").Select() | Out-Null
